$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Standardize port (POL/POD) names to Title Case across the data table.
# Column A = POL (origin port), Column B = POD (destination port).
# Column C (vehicle type) and Column D (price) are left unchanged.

$pols = @("New York", "Savannah", "Miami", "Houston", "Indianapolis", "Los Angeles", "San Francisco")
$pods = @("Rotterdam", "Varna")
$vehicleCount = 5

# First pass: write column A (POL) for every row so the new POL strings
# are registered in the shared string table before any POD strings.
$row = 2
foreach ($pod in $pods) {
    foreach ($pol in $pols) {
        for ($i = 0; $i -lt $vehicleCount; $i++) {
            $ws.Cells.Item($row, 1).Value = $pol
            $row = $row + 1
        }
    }
}

# Second pass: write column B (POD) for every row.
$row = 2
foreach ($pod in $pods) {
    foreach ($pol in $pols) {
        for ($i = 0; $i -lt $vehicleCount; $i++) {
            $ws.Cells.Item($row, 2).Value = $pod
            $row = $row + 1
        }
    }
}

# Update the active cell selection as recorded in the saved workbook.
$ws.Range("F11").Select()
